$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of credential data
$ws.Range("A3").Value = "Driver"
$ws.Range("A4").Value = "Customer"

$ws.Range("B3").Value = "driver423"

$ws.Range("B4").Value = 4565678899
$ws.Range("B4").NumberFormat = "#,##0"

# Column B should fit the widened content
$ws.Columns.Item(2).AutoFit() | Out-Null

# Update selection to match the final edited cell
$ws.Range("B4").Select() | Out-Null
